$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets("Metadata")
$meta.Range("A2").Value = "05 Nov 2025, 11:05 AM"

# --- 2. Stock List sheet: a new row is inserted at the top of the data
#         (row 2), pushing every existing data row down by one and
#         dropping the last row (76) off the bottom. ---
$ws = $wb.Worksheets("Stock List")

$lastRow = 76

# Snapshot existing values for rows 2..76 before we start overwriting them.
$B = @{}
$C = @{}
$D = @{}
$E = @{}
$H = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $B[$r] = $ws.Cells.Item($r, 2).Value2
    $C[$r] = $ws.Cells.Item($r, 3).Value2
    $D[$r] = $ws.Cells.Item($r, 4).Value2
    $E[$r] = $ws.Cells.Item($r, 5).Value2
    $H[$r] = $ws.Cells.Item($r, 8).Value2
}

# Shift rows 2..75 down into rows 3..76 (iterate from the bottom up so we
# never overwrite a value before it has been read/shifted).
for ($r = $lastRow; $r -ge 3; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 2).Value = $B[$src]
    $ws.Cells.Item($r, 3).Value = $C[$src]
    $ws.Cells.Item($r, 4).Value = $D[$src]
    $ws.Cells.Item($r, 5).Value = $E[$src]
    $ws.Cells.Item($r, 8).Value = $H[$src]
}

# New entry lands in row 2.
$ws.Cells.Item(2, 2).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 3).Value = "CAPTRU-RE1"
$ws.Cells.Item(2, 4).Value = 5.67
$ws.Cells.Item(2, 5).Value = -11.9565
$ws.Cells.Item(2, 8).Value = 0
